$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.332.11"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "2.468.33"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.63%  "
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.62"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.42%  "
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "2.849.55"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "2.454.86"
$ws.Range("E15").Value = "  -5.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.67"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.48%  "
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "41.311.28"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("E19").Value = "  -6.44%  "
$ws.Range("D20").Value = "0.0₃0921"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("E25").Value = "  -5.03%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.34%  "
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -7.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.89"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.60"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("E33").Value = "  -6.40%  "
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("E35").Value = "  -4.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("E37").Value = "  -6.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.88"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.23%  "
$ws.Range("E39").Value = "  -7.01%  "
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.30"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "1.987.09"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.73"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "69.95"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("E51").Value = "  -6.06%  "
